$wb = $excel.ActiveWorkbook

# 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1586
$ws.Range("F5").Value = 9062
$ws.Range("F10").Value = 591
$ws.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202405/rNbVBwPd1716434800421.jpeg"
$ws.Range("F17").Value = 1491
$ws.Range("F18").Value = 1320
$ws.Range("F21").Value = 1369
$ws.Range("F22").Value = 80
$ws.Range("F23").Value = 230
$ws.Range("F25").Value = 87
$ws.Range("F26").Value = 61
$ws.Range("F27").Value = 64
$ws.Range("F28").Value = 302
$ws.Range("F29").Value = 302
$ws.Range("F30").Value = 1067
$ws.Range("F34").Value = 194
$ws.Range("F42").Value = 31
$ws.Range("F46").Value = 209

# 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 161
$ws.Range("F16").Value = 669
$ws.Range("F23").Value = 930
$ws.Range("F24").Value = 16
$ws.Range("F29").Value = 217

# 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 144
$ws.Range("F7").Value = 2059
$ws.Range("F8").Value = 3092

# 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1586
$ws.Range("F5").Value = 9062
$ws.Range("F7").Value = 144
$ws.Range("F11").Value = 2059
$ws.Range("F12").Value = 3092
$ws.Range("F17").Value = 161
$ws.Range("F18").Value = 591
$ws.Range("I22").Value = "//i2.hdslb.com/bfs/openplatform/202405/rNbVBwPd1716434800421.jpeg"
$ws.Range("F23").Value = 1491
$ws.Range("F25").Value = 1369
$ws.Range("F26").Value = 230
$ws.Range("F28").Value = 87
$ws.Range("F29").Value = 302
$ws.Range("F30").Value = 302
$ws.Range("F32").Value = 930
$ws.Range("F34").Value = 16
$ws.Range("F40").Value = 217
$ws.Range("F47").Value = 209
